# Add a new "features" field row to the MultiPulse feature-extraction
# schema table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 23: Field="features", Type="cell", Description="features"
$ws.Cells.Item(23, 1).Value = "features"
$ws.Cells.Item(23, 2).Value = "cell"
$ws.Cells.Item(23, 3).Value = "features"

# Move/save the selection to C24, matching the author's saved cursor position.
$ws.Range("C24").Select()
